# Spring 23 week 13 inputs — append 14 new matchup rows (877-890) to the
# "Nine" sheet, columns A:D = Player_1, Points_1, Player_2, Points_2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$data = @(
    @(6,17,5,3),
    @(6,2,5,18),
    @(3,17,4,3),
    @(5,4,6,16),
    @(4,4,5,16),
    @(5,12,7,8),
    @(2,15,3,5),
    @(5,15,4,5),
    @(7,6,6,14),
    @(2,6,4,14),
    @(5,6,6,14),
    @(5,14,4,6),
    @(2,15,4,5),
    @(7,17,5,3)
)

$startRow = 877
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
    $ws.Range("D$row").Value = $data[$i][3]
}

# Match the post-edit view: scrolled down so the new rows are visible, with
# the next empty row selected ready for further input.
$nextRow = $endRow + 1
$ws.Application.ActiveWindow.ScrollRow = $startRow - 5
$ws.Range("A$nextRow").Select()
